# Applies updated Price (D) and Volume(1h) (E) values for the symbol list refresh
# commit: "Updated symbol list on Sun Jan  8 18:35:25 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'265.76"
$ws.Range("E2").Value = "'1.61%"

$ws.Range("D3").Value = "'26.70"
$ws.Range("E3").Value = "'-1.68%"

$ws.Range("D4").Value = "'4.704"
$ws.Range("E4").Value = "'0.12%"

$ws.Range("D5").Value = "'0.06083"
$ws.Range("E5").Value = "'-1.73%"

$ws.Range("D6").Value = "'6.695"
$ws.Range("E6").Value = "'-0.31%"

$ws.Range("D7").Value = "'0.8502"
$ws.Range("E7").Value = "'0.02%"

$ws.Range("E8").Value = "'-1.13%"

$ws.Range("D9").Value = "'0.1407"
$ws.Range("E9").Value = "'0.05%"

$ws.Range("D10").Value = "'0.04865"
$ws.Range("E10").Value = "'4.95%"

$ws.Range("D11").Value = "'0.07109"
$ws.Range("E11").Value = "'0.32%"

$ws.Range("D12").Value = "'0.03164"
$ws.Range("E12").Value = "'1.33%"

$ws.Range("D13").Value = "'0.09027"
$ws.Range("E13").Value = "'-0.09%"

$ws.Range("D14").Value = "'0.001532"
$ws.Range("E14").Value = "'-0.03%"

$ws.Range("D15").Value = "'0.0006064"
$ws.Range("E15").Value = "'-1.59%"

$ws.Range("D16").Value = "'0.005950"
$ws.Range("E16").Value = "'-3.40%"

$ws.Range("D17").Value = "'3.457"
$ws.Range("E17").Value = "'-0.04%"

$ws.Range("D18").Value = "'3.172"
$ws.Range("E18").Value = "'0.11%"

$ws.Range("D19").Value = "'2.277"
$ws.Range("E19").Value = "'3.71%"

$ws.Range("D20").Value = "'0.3087"
$ws.Range("E20").Value = "'0.51%"

$ws.Range("D21").Value = "'0.1301"
$ws.Range("E21").Value = "'0.12%"

$ws.Range("D22").Value = "'4.089"
$ws.Range("E22").Value = "'-0.16%"

$ws.Range("D23").Value = "'0.04226"
$ws.Range("E23").Value = "'-0.32%"

$ws.Range("E24").Value = "'-3.05%"

$ws.Range("D25").Value = "'0.004137"
$ws.Range("E25").Value = "'8.81%"

$ws.Range("E26").Value = "'0.04%"

$ws.Range("E27").Value = "'5.10%"

$ws.Range("D40").Value = "'0.03925"
$ws.Range("E40").Value = "'-0.26%"

$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'0.29%"

$ws.Range("D42").Value = "'0.004176"
$ws.Range("E42").Value = "'1.28%"

$ws.Range("E43").Value = "'-3.32%"

$ws.Range("D44").Value = "'0.01156"
$ws.Range("E44").Value = "'-16.90%"

$ws.Range("D45").Value = "'0.00005118"
$ws.Range("E45").Value = "'-0.24%"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.06%"

$ws.Range("D48").Value = "'0.1429"
$ws.Range("E48").Value = "'-14.80%"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.06%"

$ws.Range("E50").Value = "'0.06%"

